$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing changed cells (rows 4, 8, 11) per ddof=1 recalculation ---
$ws.Range("G4").Value = 140
$ws.Range("H4").Value = 1149
$ws.Range("I4").Value = 961
$ws.Range("J4").Value = -16.94
$ws.Range("K4").Value = 51.55
$ws.Range("J8").Value = -12.58
$ws.Range("K8").Value = 140.33
$ws.Range("G11").Value = 53
$ws.Range("H11").Value = 754
$ws.Range("I11").Value = 693
$ws.Range("J11").Value = -4.18
$ws.Range("K11").Value = 91.82

# --- Step 2: set row heights for newly introduced rows (17 existing rows 17-22 already have ht=15; 23-79 are new) ---
$ws.Range("A23:A79").EntireRow.RowHeight = 15

# --- Step 3: apply cell formatting matching style used for J/K columns (font Arial 12, centered) to the new data block ---
$ws.Range("J4").Copy()
$ws.Range("A17:K24").PasteSpecial(-4122)
$ws.Range("A25:S33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: fill new row data (column order per row matches original authoring order: A, C, D, E, F, G, H, I, J, K) ---
# Row 17
$ws.Range("A17").Value = "FatimaAparecida"
$ws.Range("C17").Value = 451
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 67
$ws.Range("F17").Value = 374
$ws.Range("G17").Value = 67
$ws.Range("H17").Value = 883
$ws.Range("I17").Value = 782
$ws.Range("J17").Value = -17.01
$ws.Range("K17").Value = 37.28
# Row 18
$ws.Range("A18").Value = "MarliRosa"
$ws.Range("C18").Value = 493
$ws.Range("D18").Value = 905
$ws.Range("E18").Value = 46
$ws.Range("F18").Value = 384
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = 938
$ws.Range("I18").Value = 840
$ws.Range("J18").Value = -20.02
$ws.Range("K18").Value = 41.47
# Row 19
$ws.Range("A19").Value = "JoseMacedo"
$ws.Range("C19").Value = "SEMTEMPOS"
# Row 20
$ws.Range("A20").Value = "FabricioLucena"
$ws.Range("C20").Value = 397
$ws.Range("D20").Value = 877
$ws.Range("E20").Value = 65
$ws.Range("F20").Value = 326
$ws.Range("G20").Value = 48
$ws.Range("H20").Value = 942
$ws.Range("I20").Value = 835
$ws.Range("J20").Value = -19.29
$ws.Range("K20").Value = 37.97
# Row 21
$ws.Range("A21").Value = "Neusa"
$ws.Range("C21").Value = 515
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 31
$ws.Range("F21").Value = 385
$ws.Range("G21").Value = 119
$ws.Range("H21").Value = 1216
$ws.Range("I21").Value = 1113
$ws.Range("J21").Value = -15.15
$ws.Range("K21").Value = 69.8
# Row 22
$ws.Range("A22").Value = "Eliene"
$ws.Range("C22").Value = 484
$ws.Range("D22").Value = 1071
$ws.Range("E22").Value = 78
$ws.Range("F22").Value = 396
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = 1014
$ws.Range("I22").Value = 784
$ws.Range("J22").Value = -16.42
$ws.Range("K22").Value = 89.34
# Row 23
$ws.Range("A23").Value = "Adauto"
$ws.Range("C23").Value = 561
$ws.Range("D23").Value = 45
$ws.Range("E23").Value = 68
$ws.Range("F23").Value = 412
$ws.Range("G23").Value = 119
$ws.Range("H23").Value = 1255
$ws.Range("I23").Value = 1099
$ws.Range("J23").Value = -18.99
$ws.Range("K23").Value = 49.95
# Row 24
$ws.Range("A24").Value = "JoseSilva"
$ws.Range("C24").Value = 453
$ws.Range("D24").Value = 924
$ws.Range("E24").Value = 49
$ws.Range("F24").Value = 355
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = 1153
$ws.Range("I24").Value = 950
$ws.Range("J24").Value = -16.81
$ws.Range("K24").Value = 51.55
# Row 25
$ws.Range("A25").Value = "Diogenes"
$ws.Range("C25").Value = 399
$ws.Range("D25").Value = 796
$ws.Range("E25").Value = 74
$ws.Range("F25").Value = 354
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 1162
$ws.Range("I25").Value = 1043
$ws.Range("J25").Value = -10.3
$ws.Range("K25").Value = 96.77
# Row 26
$ws.Range("A26").Value = "MariaSilva"
$ws.Range("C26").Value = 548
$ws.Range("D26").Value = 55
$ws.Range("E26").Value = 78
$ws.Range("F26").Value = 409
$ws.Range("G26").Value = 115
$ws.Range("H26").Value = 1280
$ws.Range("I26").Value = 1194
$ws.Range("J26").Value = -19.84
$ws.Range("K26").Value = 36.86
# Row 27
$ws.Range("A27").Value = "MarcoAntonio"
$ws.Range("C27").Value = 730
$ws.Range("D27").Value = 1011
$ws.Range("E27").Value = 1085
$ws.Range("F27").Value = 1333
# Row 28
$ws.Range("A28").Value = "Marlene"
$ws.Range("C28").Value = 493
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 63
$ws.Range("F28").Value = 425
$ws.Range("G28").Value = 35
$ws.Range("H28").Value = 939
$ws.Range("I28").Value = 851
$ws.Range("J28").Value = -20.9
$ws.Range("K28").Value = 132.72
# Row 29
$ws.Range("A29").Value = "Valdir"
$ws.Range("C29").Value = 422
$ws.Range("D29").Value = 840
$ws.Range("E29").Value = 38
$ws.Range("F29").Value = 331
$ws.Range("G29").Value = 101
$ws.Range("H29").Value = 1019
$ws.Range("I29").Value = 954
$ws.Range("J29").Value = -15.04
$ws.Range("K29").Value = 82.59
# Row 30
$ws.Range("A30").Value = "PedroFernandes"
$ws.Range("C30").Value = 434
$ws.Range("D30").Value = 976
$ws.Range("E30").Value = 40
$ws.Range("F30").Value = 341
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 1074
$ws.Range("I30").Value = 871
$ws.Range("J30").Value = -19.96
$ws.Range("K30").Value = 101.56
# Row 31
$ws.Range("A31").Value = "MariaNogueira"
$ws.Range("C31").Value = 351
$ws.Range("D31").Value = 26
$ws.Range("E31").Value = 74
$ws.Range("F31").Value = 384
# Row 32
$ws.Range("A32").Value = "Patricia"
$ws.Range("C32").Value = 445
$ws.Range("D32").Value = 817
$ws.Range("E32").Value = 60
$ws.Range("F32").Value = 399
$ws.Range("G32").Value = 29
$ws.Range("H32").Value = 837
$ws.Range("I32").Value = 761
$ws.Range("J32").Value = -20.42
$ws.Range("K32").Value = 53
# Row 33
$ws.Range("A33").Value = "Roseli"
$ws.Range("C33").Value = 443
$ws.Range("D33").Value = 783
$ws.Range("E33").Value = 61
$ws.Range("F33").Value = 359
$ws.Range("G33").Value = 20
$ws.Range("H33").Value = 698
$ws.Range("I33").Value = 610
$ws.Range("J33").Value = -17.96
$ws.Range("K33").Value = 150.33

# --- Step 5: set active selection to match authors final cursor position ---
$ws.Range("B25").Select()
